$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.020378666666666
$ws.Range("H2").Value = 6.061135999999999
$ws.Range("I2").Value = 0.2600911804892308
$ws.Range("J2").Value = 0.2600911804892308
$ws.Range("Q2").Value = 0.1661492069511111
$ws.Range("R2").Value = 1.49534286256
$ws.Range("S2").Value = 0.2600911804892308
$ws.Range("T2").Value = 0.2600911804892308

# Row 3
$ws.Range("I3").Value = 0.4664203043534923
$ws.Range("J3").Value = 0.4664203043534922
$ws.Range("R3").Value = 2.68159140097
$ws.Range("S3").Value = 0.4664203043534923
$ws.Range("T3").Value = 0.4664203043534922

# Row 4
$ws.Range("G4").Value = 2.124448666666666
$ws.Range("H4").Value = 6.373346
$ws.Range("I4").Value = 0.273488515157277
$ws.Range("J4").Value = 0.273488515157277
$ws.Range("Q4").Value = 0.1747075768511111
$ws.Range("R4").Value = 1.57236819166
$ws.Range("S4").Value = 0.273488515157277
$ws.Range("T4").Value = 0.273488515157277
